$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the anchor heading "Philippians 1:1" (Heading 2). The two
# paragraphs that need new content sit just before it:
#   5 paragraphs back -> the (currently empty) book-abbreviation
#                         Heading 2 paragraph -> gets "PHP"
#   4 paragraphs back -> the (currently empty) Normal paragraph that
#                         holds the italic verse-reference list
# ------------------------------------------------------------------
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.ParagraphStyle.NameLocal -eq "Heading 2" -and $p.Range.Text -eq "Philippians 1:1`r") {
        $anchor = $p
        break
    }
}

$bookAbbrevPara = $anchor.Previous(5)
$verseListPara = $anchor.Previous(4)

# --- Paragraph: book abbreviation heading ("PHP") ---
$bookAbbrevPara.Range.Text = "PHP"

# --- Paragraph: italic list of every verse reference in the book ---
$verses = @()
for ($v = 1; $v -le 30; $v++) { $verses += "Philippians 1:$v" }
for ($v = 1; $v -le 30; $v++) { $verses += "Philippians 2:$v" }
for ($v = 1; $v -le 21; $v++) { $verses += "Philippians 3:$v" }
for ($v = 1; $v -le 23; $v++) { $verses += "Philippians 4:$v" }
$verseText = [string]::Join(", ", $verses)

$verseListPara.Range.Text = $verseText
$verseListStart = $verseListPara.Range.Start
$verseTextRange = $d.Range($verseListStart, $verseListStart + $verseText.Length)
$verseTextRange.Font.Italic = 1

Write-Host "Book abbreviation paragraph now: [" $bookAbbrevPara.Range.Text "]"
Write-Host "Verse list paragraph length: " $verseListPara.Range.Text.Length
